$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data (dates as Excel serial numbers, matching existing style)
$ws.Range("A119").Value = 45875
$ws.Range("B119").Value = 1.88

$ws.Range("A120").Value = 45889
$ws.Range("B120").Value = 1.94

# Copy number formats from the row above to keep formatting consistent
$ws.Range("A118").Copy() | Out-Null
$ws.Range("A119:A120").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B118").Copy() | Out-Null
$ws.Range("B119:B120").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the selected cell to match the new final cell
$ws.Range("B120").Select() | Out-Null
